$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated voltage magnitude (vm_pu) results for the 380 kV case.
# Each inner array holds: row number, followed by new values for columns B,C,D,E,F,I,J,K,L,M,N
# (column G stays fixed at 1 and column H is unused, so both are skipped).
$colNumbers = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(2,1.02,1.063284187537301,1.073204712324558,1.069318210248063,1.081502354518122,1.044216025439591,1.068250355646021,1.075897192068898,1.072021038038645,1.084173010107884,1.026377243162983),
    @(3,1.02,1.064641936880491,1.074506230319137,1.070546795250658,1.082866962142897,1.044501256903282,1.069260929125355,1.077014408137109,1.073064747482746,1.085354752786499,1.026719167937184),
    @(4,1.02,1.065519469045442,1.075347709455992,1.071341195705979,1.0837494935339,1.044683387587444,1.069913326073978,1.077736077761744,1.073738953239911,1.086118413953683,1.026939760598712),
    @(5,1.02,1.065888143057409,1.075701305561853,1.07167502686547,1.08412040299614,1.044759373705639,1.070187234766378,1.078039173018811,1.074022118349826,1.086439219678486,1.027032341736111),
    @(6,1.02,1.065950031066162,1.075760666520879,1.071731070772967,1.084182674172967,1.044772098027912,1.070233204274237,1.078090046858197,1.074069647210538,1.086493070542525,1.027047877394065),
    @(7,1.02,1.065524396225398,1.075352434855638,1.07134565690002,1.083754450059197,1.044684405201989,1.069916987462314,1.077740128889328,1.073742737967207,1.086122701504399,1.026940998284313),
    @(8,1.02,1.063743259565192,1.073644711398267,1.06973353750939,1.08196362769562,1.044312925080776,1.068592198444845,1.076275019916089,1.072374003163168,1.084572595714006,1.026492934528058),
    @(9,1.02,1.060596646342874,1.070630042967198,1.066888222063316,1.078804275694183,1.043639654337771,1.066246041847401,1.073683642627361,1.069953224826274,1.081833264395395,1.025698326864108),
    @(10,1.02,1.058493238564244,1.068616383193573,1.064988080701654,1.076695325438751,1.043178192236346,1.0646738780158,1.071949356602857,1.06833320800119,1.080001559191932,1.025165125856804),
    @(11,1.02,1.057581034672417,1.067743475378615,1.064164478416264,1.075781430062544,1.042975368165945,1.063991159624976,1.071196756545122,1.067630218500516,1.079207062480211,1.024933409170432),
    @(12,1.02,1.057241983638558,1.067419086883225,1.063858427324204,1.075441857925989,1.042899577082566,1.063737269212264,1.070916956522845,1.067368866071941,1.078911743074227,1.024847212301132),
    @(13,1.02,1.057314721183283,1.067488676251813,1.063924082204167,1.075514702386303,1.042915855055983,1.063791743132013,1.070976985994271,1.067424937560818,1.078975099567173,1.024865707587535),
    @(14,1.02,1.057553013084439,1.067716664444251,1.064139182783109,1.07575336320904,1.042969112503844,1.063970179073567,1.071173633328802,1.06760861975778,1.079182655552495,1.024926286707582),
    @(15,1.02,1.057699803554642,1.067857115262957,1.064271696195397,1.075900395210318,1.043001866083989,1.064080079691356,1.071294760879998,1.067721761722401,1.079310509963722,1.024963594641187),
    @(16,1.02,1.05855374827252,1.06867429418665,1.065042722671247,1.076755962395691,1.043191589518945,1.064719146147819,1.071999269234433,1.068379830934655,1.080054258390303,1.025180486390093),
    @(17,1.02,1.059089022753854,1.069186623468274,1.065526142534002,1.07729244481024,1.043309791806242,1.065119487960754,1.072440746052162,1.068792213057546,1.080520426069813,1.025316311882192),
    @(18,1.02,1.059401103196398,1.069485362233437,1.065808033495928,1.077605298123436,1.043378447066805,1.065352811609086,1.072698093783721,1.06903260298764,1.080792203215479,1.02539545590498),
    @(19,1.02,1.059507491638451,1.069587208581875,1.065904137638363,1.077711961657092,1.043401807563077,1.065432337006331,1.072785815993121,1.06911454516003,1.080884850176812,1.025422428321437),
    @(20,1.02,1.059031606986444,1.069131665166004,1.065474284428559,1.077234892372505,1.043297139839186,1.065076554646342,1.072393396164427,1.068747983445311,1.080470424227165,1.025301747447481),
    @(21,1.02,1.057482848117184,1.06764953182674,1.064075844604693,1.075683086638539,1.042953442033606,1.06391764240407,1.071115732547459,1.067554536331069,1.079121541262314,1.024908451181478),
    @(22,1.02,1.056507818591015,1.066716776041857,1.06319584430403,1.074706761970877,1.042734722898241,1.063187259912786,1.070310961346156,1.06680283139366,1.078272239759763,1.024660434968086),
    @(23,1.02,1.057024821750628,1.067211332416021,1.063662421156557,1.075224392636663,1.042850919083668,1.063574614587824,1.070737724908126,1.067201452425957,1.078722586205044,1.02479198315255),
    @(24,1.02,1.059057551120537,1.069156498748108,1.065497717120169,1.077260898051076,1.043302857614918,1.065095954946737,1.072414792029234,1.068767969352888,1.0804930183108,1.025308328737871),
    @(25,1.02,1.06141109986818,1.071410074633784,1.067624364784374,1.079621505451325,1.043815930061817,1.066853986290042,1.074354742064357,1.070580126900414,1.082542398459031,1.025904357926871)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $colNumbers.Length; $i++) {
        $ws.Cells.Item($r, $colNumbers[$i]).Value = $rowData[$i + 1]
    }
}
